# Insert a new data row at row 81 (pushes existing rows 81-120 down to 82-121)
# and populate it with the new weekly price record for Cebollín.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Insert()

$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = 45097
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 100112037
$ws.Range("G81").Value = "Cebollín"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 230
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = 4326
$ws.Range("N81").Value = "$/paquete 36 unidades"
$ws.Range("O81").Value = "Región Metropolitana"
$ws.Range("P81").Value = 120
$ws.Range("Q81").Value = 36
$ws.Range("R81").Value = "Hortaliza"
